$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New JLCPCB submission: C16 (10uF tantalum cap) LCSC part number corrected
# from C307331 to C110055.
$ws.Range("D5").Value = "C110055"

# Widen the Footprint column so the long footprint names are fully visible.
$ws.Columns("C").ColumnWidth = 57
